# Update the NATMI TPM-derived LR-pair stats for Il17f-Il17ra with new TPM values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Sending cluster=MuSCs, Ligand=Il17f, Receptor=Il17ra, Target cluster=ECs
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.05839233333333333
$ws.Range("H2").Value = 0.175177
$ws.Range("M2").Value = 0.1801153333333333
$ws.Range("N2").Value = 0.540346
$ws.Range("O2").Value = 0.01663333613045927
$ws.Range("P2").Value = 0.01663333613045927
$ws.Range("Q2").Value = 0.01051735458244444
$ws.Range("R2").Value = 0.094656191242
$ws.Range("S2").Value = 0.01663333613045927
$ws.Range("T2").Value = 0.01663333613045927

# Row 3: Sending cluster=MuSCs, Ligand=Il17f, Receptor=Il17ra, Target cluster=FAPs
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.05839233333333333
$ws.Range("H3").Value = 0.175177
$ws.Range("O3").Value = 0.407089716880577
$ws.Range("P3").Value = 0.407089716880577
$ws.Range("Q3").Value = 0.2574051811205555
$ws.Range("R3").Value = 2.316646630085
$ws.Range("S3").Value = 0.407089716880577
$ws.Range("T3").Value = 0.407089716880577

# Row 4: Sending cluster=MuSCs, Ligand=Il17f, Receptor=Il17ra, Target cluster=MuSCs
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.05839233333333333
$ws.Range("H4").Value = 0.175177
$ws.Range("M4").Value = 6.240258333333333
$ws.Range("N4").Value = 18.720775
$ws.Range("O4").Value = 0.5762769469889637
$ws.Range("P4").Value = 0.5762769469889637
$ws.Range("Q4").Value = 0.3643832446861111
$ws.Range("R4").Value = 3.279449202175
$ws.Range("S4").Value = 0.5762769469889637
$ws.Range("T4").Value = 0.5762769469889637
